$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C17").Value = "Travail à réaliser : "
$ws.Range("C18").Value = "Réaliser les routes de l'API Fastify."
$ws.Range("C19").Value = "Faire les tables avec MIKRO ORM."
$ws.Range("C20").Value = "Réaliser le protocole en JAVA."

$ws.Range("C17:C20").Font.Bold = $true
$ws.Range("C17:C20").Font.Size = 20

$ws.Range("C21").Select()
